$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header column G1
$ws.Range("G1").Value = "estado pedido"

# Update row 2 values
$ws.Range("A2").Value = "prueba 1"
$ws.Range("B2").Value = "40 kg"
$ws.Range("C2").Value = "basico"
$ws.Range("D2").Value = "arroz,  papa,  yuca"
$ws.Range("E2").Value = "alimentos"
$ws.Range("F2").Value = "20x20x20"
$ws.Range("G2").Value = "pendiente"

# Delete row 3 entirely (removing the second "prueba" entry)
$ws.Rows.Item(3).Delete()
